$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (id=4, class=1) gets new course codes.
# Shared-string append order observed in the target file is F5, E5, D5, C5,
# so write the cells in that same order to reproduce it.
$ws.Range("F5").Value = "0DKoqD2pKiDU5LniiSdz"
$ws.Range("E5").Value = "WsiJrcsuBgOwJoji2lx3"
$ws.Range("D5").Value = "44qXU7UPPjExGS7x1Wpg"
$ws.Range("C5").Value = "GlAcMvhEhyOA7GMV64K1"

# Row 11 (id=10, class=2) gets new course codes, written left-to-right.
$ws.Range("C11").Value = "82it5OgUXmN9955YfAw0"
$ws.Range("D11").Value = "mcIUEEaoTvKKrD88MXDJ"
$ws.Range("E11").Value = "RSYlMwCHc6qAHUHbRERp"
$ws.Range("F11").Value = "C5n7N4ywR1zCZc74ktyc"
